# AVA sheet: collapse the TimeSlice/PRE/DMD availability rows into a single
# PSET_PN / "*" row (commodity TSL cases).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVA")

# Header row (row 4): rename the first column heading "TimeSlice" -> "PSET_PN"
$ws.Range("B4").Value = "PSET_PN"

# Remove the old row 5 (the "PRE" marker row, F5=0,G5=1,H5=1); this shifts
# the old row 6 (the "DMD" marker row, F6=1, I6:P6=0) up to become row 5.
$ws.Rows.Item(5).Delete()

# Turn what is now the single remaining data row into the "*" (all/default)
# row: give it a row label and drop the now-unused trailing marker value.
$ws.Range("B5").Value = "*"
$ws.Range("Q5").ClearContents()

$ws.Range("B6").Select()
